# Generate Report for Archive
# - Update status text "Ready for handoff" -> "In Translation"
# - Shrink the (now narrower) Status/zh-cn/de-de columns to match the new content width

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# Target stored column width (OOXML "width" units) that the columns should end up at.
# Derived from: stored_width = ColumnWidth_input + 0.8333333333333334 (rounded to the
# nearest 1/6 by this host), so we back-solve for the ColumnWidth value to assign.
$targetStoredWidth = 13.4101845877511
$columnWidthToSet = $targetStoredWidth - 0.8333333333333334

# --- Overview sheet: columns E (zh-cn) and F (de-de) both hold the status text ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = $columnWidthToSet
$wsOverview.Columns.Item(6).ColumnWidth = $columnWidthToSet

# --- zh-cn sheet: column C (Status) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = $columnWidthToSet

# --- de-de sheet: column C (Status) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = $columnWidthToSet
